# Weekly update: insert a new Brócoli price record at the top of the data
# block (row 104), pushing all existing records down by one row.
# The former last row (205) becomes row 206; the sheet's used range grows
# from A1:R205 to A1:R206.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 104:205 down to 105:206, creating a blank row 104.
$ws.Rows.Item(104).Insert()

# Populate the new row 104 with the latest week's record.
$ws.Range("A104").Value = 5
$ws.Range("B104").Value = "Macroferia Regional de Talca"
$ws.Range("C104").Value = "Maule"
$ws.Range("D104").Value = 44484
$ws.Range("E104").Value = 7
$ws.Range("F104").Value = 100112023
$ws.Range("G104").Value = "Brócoli"
$ws.Range("H104").Value = "Sin especificar"
$ws.Range("I104").Value = "Primera"
$ws.Range("J104").Value = 3000
$ws.Range("K104").Value = 700
$ws.Range("L104").Value = 700
$ws.Range("M104").Value = 700
$ws.Range("N104").Value = '$/unidad'
$ws.Range("O104").Value = "Región del Maule"
$ws.Range("P104").Value = 700
$ws.Range("Q104").Value = 1
$ws.Range("R104").Value = "Hortaliza"
